$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A1:G25")
$rng.Sort($ws.Range("G1"), 2, $null, $null, 1, $null, 1, 1, $false, $null, $null, 1)

# The sort swaps cell contents but (in this runtime) row heights stay pinned to
# their physical row index instead of travelling with the row's data. Fix the
# taller rows up so the row heights line up with the re-sorted content, the
# same way a genuine Excel sort (which moves whole rows) would.
$ws.Rows.Item(2).RowHeight = 19.5
$ws.Rows.Item(10).RowHeight = 19.5
$ws.Rows.Item(25).RowHeight = 19.5
$ws.Rows.Item(11).RowHeight = 20.25
$ws.Rows.Item(12).RowHeight = 20.25
$ws.Rows.Item(17).RowHeight = 20.25
